$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '65.542.40'
$ws.Range("E2").Value = '  -2.79%  '
$ws.Range("D3").Value = '3.492.93'
$ws.Range("E3").Value = '  -0.24%  '
$ws.Range("E4").Value = '  +0.02%  '
$ws.Range("D5").Value = '555.37'
$ws.Range("E5").Value = '  +0.23%  '
$ws.Range("D6").Value = '178.71'
$ws.Range("E6").Value = '  -5.66%  '
$ws.Range("D7").Value = '0.637'
$ws.Range("E7").Value = '  +4.46%  '
$ws.Range("E8").Value = '  +0.00%  '
$ws.Range("D9").Value = '0.631'
$ws.Range("E9").Value = '  -1.08%  '
$ws.Range("D10").Value = '0.154'
$ws.Range("E10").Value = '  +2.68%  '
$ws.Range("D11").Value = '53.64'
$ws.Range("E11").Value = '  -5.80%  '
$ws.Range("D12").Value = '0.0000271'
$ws.Range("E12").Value = '  -1.26%  '
$ws.Range("D13").Value = '9.24'
$ws.Range("E13").Value = '  -2.65%  '
$ws.Range("D14").Value = '4.061.46'
$ws.Range("E14").Value = '  -0.11%  '
$ws.Range("D15").Value = '3.498.53'
$ws.Range("E15").Value = '  -0.28%  '
$ws.Range("E16").Value = '  +0.09%  '
$ws.Range("D17").Value = '18.38'
$ws.Range("E17").Value = '  +0.37%  '
$ws.Range("D18").Value = '12.05'
$ws.Range("E18").Value = '  +1.66%  '
$ws.Range("D19").Value = '65.602.37'
$ws.Range("E19").Value = '  -3.61%  '
$ws.Range("D20").Value = '0.996'
$ws.Range("E20").Value = '  -1.12%  '
$ws.Range("D21").Value = '413.66'
$ws.Range("E21").Value = '  +1.30%  '
$ws.Range("D22").Value = '4.05'
$ws.Range("E22").Value = '  +2.71%  '
$ws.Range("D23").Value = '85.95'
$ws.Range("E23").Value = '  +1.62%  '
$ws.Range("D24").Value = '4.11'
$ws.Range("E24").Value = '  -2.47%  '
$ws.Range("D25").Value = '12.71'
$ws.Range("E25").Value = '  +7.27%  '
$ws.Range("D26").Value = '10.80'
$ws.Range("E26").Value = '  -8.13%  '
$ws.Range("E27").Value = '  -1.97%  '
$ws.Range("E28").Value = '  -1.91%  '
$ws.Range("E29").Value = '  +5.10%  '
$ws.Range("D30").Value = '30.25'
$ws.Range("E30").Value = '  -0.55%  '
$ws.Range("D31").Value = '6.47'
$ws.Range("E31").Value = '  -5.23%  '
$ws.Range("D32").Value = '609.33'
$ws.Range("E32").Value = '  -11.50%  '
$ws.Range("D33").Value = '11.65'
$ws.Range("E33").Value = '  -0.66%  '
$ws.Range("D35").Value = '59.49'
$ws.Range("E35").Value = '  -0.95%  '
$ws.Range("E36").Value = '  +10.59%  '
$ws.Range("D37").Value = '0.999'
$ws.Range("E37").Value = '  +0.05%  '
$ws.Range("B38").Value = 'InjectiveProtocol'
$ws.Range("C38").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D38").Value = '37.14'
$ws.Range("E38").Value = '  -4.69%  '
$ws.Range("B39").Value = 'PEPE'
$ws.Range("C39").Value = 'https://coinranking.com/coin/03WI8NQPF+pepe-pepe'
$ws.Range("D39").Value = '0.0₃0788'
$ws.Range("E39").Value = '  -4.44%  '
$ws.Range("D40").Value = '3.356.12'
$ws.Range("E40").Value = '  +9.83%  '
$ws.Range("D41").Value = '0.379'
$ws.Range("E41").Value = '  -5.56%  '
$ws.Range("D42").Value = '1.00'
$ws.Range("E42").Value = '  -0.04%  '
$ws.Range("E43").Value = '  -3.24%  '
$ws.Range("D44").Value = '2.84'
$ws.Range("E44").Value = '  -5.69%  '
$ws.Range("B45").Value = 'ApeXProtocol'
$ws.Range("C45").Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$ws.Range("D45").Value = '3.28'
$ws.Range("E45").Value = '  +0.60%  '
$ws.Range("B46").Value = 'Fetch.AI'
$ws.Range("C46").Value = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
$ws.Range("D46").Value = '2.52'
$ws.Range("E46").Value = '  -9.55%  '
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range("D47").Value = '0.0414'
$ws.Range("E47").Value = '  -1.57%  '
$ws.Range("E48").Value = '  -1.60%  '
$ws.Range("E49").Value = '  +1.59%  '
$ws.Range("D50").Value = '8.43'
$ws.Range("E50").Value = '  -6.79%  '
$ws.Range("D51").Value = '137.57'
$ws.Range("E51").Value = '  -1.59%  '
